$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Restructure the "Create granule level metadata..." paragraph:
#    - split the word "image" (the first occurrence, right after "each ")
#      into its own run and replace it with "level one data"
#    - move the "_GoBack" bookmark from its old home (end of the
#      "Describe how to access the database..." paragraph) to right
#      after the new "level one data" run
# ---------------------------------------------------------------------

$granuleParaText = "Create granule level metadata for each image with a link to a reduced quality image."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText.StartsWith("Create granule level metadata for each image")) {
        $target = $d.Paragraphs.Item($i)
        break
    }
}

$paraStart = $target.Range.Start
$paraEnd = $target.Range.End
# "image" is the first occurrence of that word in the paragraph, right
# after "Create granule level metadata for each "
$prefix = "Create granule level metadata for each "
$imgStart = $paraStart + $prefix.Length
$imgEnd = $imgStart + "image".Length

# Force a run boundary right before "image" using a throwaway bookmark,
# then delete the bookmark (the run split survives the delete).
$splitPoint = $d.Range($imgStart, $imgStart)
$d.Bookmarks.Add("ZZSPLITTMP", $splitPoint) | Out-Null

# Replace "image" with "level one data" (this run now sits between the
# split point above and the remainder of the sentence).
$wordRange = $d.Range($imgStart, $imgEnd)
$wordRange.Text = "level one data"
$afterLevelOneData = $imgStart + "level one data".Length

# Re-home the (unique, Word-enforced) "_GoBack" bookmark immediately
# after "level one data" -- adding it here automatically removes it
# from wherever it used to be.
$bmPoint = $d.Range($afterLevelOneData, $afterLevelOneData)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# Clean up the throwaway bookmark used only to force the run split.
$d.Bookmarks("ZZSPLITTMP").Delete()

# ---------------------------------------------------------------------
# 2) Add a new paragraph after "List everyone working on the
#    telescopes..." with the new FITS-headers sentence.
# ---------------------------------------------------------------------

$listPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText.StartsWith("List everyone working on the telescopes")) {
        $listPara = $d.Paragraphs.Item($i)
        break
    }
}

$insertionRange = $listPara.Range
$insertionRange.Collapse(0)
$insertionRange.InsertParagraphAfter()
$insertionRange.Collapse(0)
$newIndex = $listPara.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Text = "Increase amount of information in the metadata that is in the FITS headers "
